# Workbook: "Fruta, Terminal Hortofrutícola Agro Chillán - Plátano"
# Insert two new price rows at the top of the "Pintón" / "Primera Pintón"
# block (rows 564-565), pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 564 (this shifts rows 564:647 down to 566:649
# and grows the used range from A1:T647 to A1:T649).
$ws.Range("A564:A565").EntireRow.Insert()

# --- New row 564 ---
$ws.Cells.Item(564, 1).Value = 7
$ws.Cells.Item(564, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(564, 3).Value = "Ñuble"
$ws.Cells.Item(564, 4).Value = 44776
$ws.Cells.Item(564, 5).Value = 16
$ws.Cells.Item(564, 6).Value = "Fruta"
$ws.Cells.Item(564, 7).Value = 100108
$ws.Cells.Item(564, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(564, 9).Value = 100108006
$ws.Cells.Item(564, 10).Value = "Plátano"
$ws.Cells.Item(564, 11).Value = "Sin especificar"
$ws.Cells.Item(564, 12).Value = "Pintón"
$ws.Cells.Item(564, 13).Value = 80
$ws.Cells.Item(564, 14).Value = 27000
$ws.Cells.Item(564, 15).Value = 27000
$ws.Cells.Item(564, 16).Value = 27000
$ws.Cells.Item(564, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(564, 18).Value = "Ecuador"
$ws.Cells.Item(564, 19).Value = 1350
$ws.Cells.Item(564, 20).Value = 20

# --- New row 565 ---
$ws.Cells.Item(565, 1).Value = 7
$ws.Cells.Item(565, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(565, 3).Value = "Ñuble"
$ws.Cells.Item(565, 4).Value = 44776
$ws.Cells.Item(565, 5).Value = 16
$ws.Cells.Item(565, 6).Value = "Fruta"
$ws.Cells.Item(565, 7).Value = 100108
$ws.Cells.Item(565, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(565, 9).Value = 100108006
$ws.Cells.Item(565, 10).Value = "Plátano"
$ws.Cells.Item(565, 11).Value = "Sin especificar"
$ws.Cells.Item(565, 12).Value = "Primera Pintón"
$ws.Cells.Item(565, 13).Value = 160
$ws.Cells.Item(565, 14).Value = 28000
$ws.Cells.Item(565, 15).Value = 29000
$ws.Cells.Item(565, 16).Value = 28500
$ws.Cells.Item(565, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(565, 18).Value = "Ecuador"
$ws.Cells.Item(565, 19).Value = 1425
$ws.Cells.Item(565, 20).Value = 20

Write-Host "Inserted 2 new rows with Plátano Pintón / Primera Pintón price data."
